# Modularization in progress 06/19/2020
# Update the cached "datetimeFigureOut" footer field text from 2020-06-04
# to 2020-06-19 everywhere it appears: the slide master and every one of
# its slide layouts (each has exactly one Date Placeholder shape holding
# the cached date string).

$p = $ppt.ActivePresentation

$oldDate = "2020-06-04"
$newDate = "2020-06-19"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# All slide layouts that hang off the master
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Just in case any individual slide also carries its own cached date field
for ($S = 1; $S -le $p.Slides.Count; $S++) {
    $slide = $p.Slides.Item($S)
    Update-DatePlaceholder $slide.Shapes
}
